$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("L6")
$r.Value2 = $r.Value2 + "|potteryman;1"
